$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = " 09/06/2019 7:00"
$ws.Range("B4").Value = " 09/20/2019 12:00"
$ws.Range("B3").Value = " 09/12/2019 00:00"
$ws.Range("B5").Value = " 09/24/2019 00:00"
$ws.Range("B6").Value = " 09/25/2019 7:00"
$ws.Range("B7").Value = " 10/09/2019 00:00"

$ws.Range("B7").Select()
